# Generate Report for Handback
# Rename the "fr-fr" locale sheet/table/column to "zh-cn" and refresh the
# handoff/handback timestamps recorded on the report.

$wb = $excel.ActiveWorkbook

# --- Locate the two worksheets before any renames happen ---------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsLocale   = $wb.Worksheets.Item("fr-fr")

# --- 1. Update the per-file handoff / handback timestamps on the ------
#        locale sheet (columns E and H, data rows 2-5).
$wsLocale.Range("E2:E5").Value = "2016-03-11 01:02:59"
$wsLocale.Range("H2:H5").Value = "2016-03-18 09:21:04"

# --- 2. Rename the locale table's second column on the Overview sheet -
#        (table has headerRowCount=0, so there is no visible header row
#        to edit a cell in directly -- briefly turn headers on, edit the
#        header cell, then restore the original headerless layout/ref).
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ShowHeaders = $true
$wsOverview.Range("B1").Value = "zh-cn"
$loOverview.ShowHeaders = $false
$loOverview.Resize($wsOverview.Range("A1:C1"))

# --- 3. Rename the locale table itself (on the "fr-fr" sheet). --------
$loLocale = $wsLocale.ListObjects.Item(1)
$loLocale.Name = "zh-cn"

# --- 4. Rename the worksheet tab last, from "fr-fr" to "zh-cn". -------
$wsLocale.Name = "zh-cn"
